# food_bout form: the date/time typed survey fields are being switched to
# plain text fields (FB_FOL_date, FB_begin_feed_time, FB_end_feed_time).

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

$survey.Range("C2").Value = "text"
$survey.Range("C4").Value = "text"
$survey.Range("C5").Value = "text"

# Move the cursor to where the edit was last made.
$survey.Range("C9").Select()
